# FEAT: improves states transition on EOF character.
# Updates the "C_EOF" (column AH) transitions of the lexer FSM table so
# that each state transitions to its proper token on end-of-file instead
# of uniformly emitting T_EOF. Cells whose new target is the T_ERROR
# state are also re-colored red (matching the existing T_ERROR styling
# used throughout the rest of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose AH (C_EOF) transition becomes T_ERROR -> also recolor red.
$errorRows = @(5, 7, 11, 14, 16, 17, 18, 26, 28, 30, 34, 35, 36, 37, 38)
foreach ($r in $errorRows) {
    $cell = $ws.Range("AH$r")
    $cell.Value = "T_ERROR"
    $cell.Font.Color = 192
}

# Rows whose AH (C_EOF) transition becomes some other, more specific token.
$otherRows = @{
    9  = "T_WORD"
    10 = "T_FILE"
    12 = "T_REFINE"
    13 = "T_ISSUE"
    19 = "T_ISSUE"
    20 = "T_INTEGER"
    21 = "T_FLOAT"
    22 = "T_FLOAT"
    33 = "T_WORD"
    39 = "T_WORD"
    40 = "T_WORD"
    41 = "T_WORD"
    42 = "T_URL"
    43 = "T_EMAIL"
}
foreach ($r in $otherRows.Keys) {
    $ws.Range("AH$r").Value = $otherRows[$r]
}
